$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

$updates = @(
    @{Row=2; B=0.00236655069326775; C=0.649198730991558},
    @{Row=3; B=0.00273222834405392; C=0.542272216357288},
    @{Row=4; B=0.00279248261528358; C=0.58431274528565},
    @{Row=5; B=0.00281627144216852; C=0.647318447403382},
    @{Row=6; B=0.00318674057032065; C=0.699363159911048},
    @{Row=7; B=0.00334827073445093; C=0.623141254306129},
    @{Row=8; B=0.00321848395386037; C=0.746260244142302},
    @{Row=9; B=0.00264247241616706; C=0.739854415825333},
    @{Row=10; B=0.00289923695065701; C=0.649918885562665},
    @{Row=11; B=0.0027343908310323; C=0.643772520705782},
    @{Row=12; B=0.00543972911234428; C=0.62065986148731},
    @{Row=13; B=0.00309196118596838; C=0.766862976337199},
    @{Row=14; B=0.00305906564565818; C=0.682143343498224},
    @{Row=15; B=0.00281137806399219; C=0.780488510496068},
    @{Row=16; B=0.00267788758712278; C=0.609383735684052},
    @{Row=17; B=0.00333296035881704; C=0.613977858195642},
    @{Row=18; B=0.00259529550156573; C=0.711548390955164},
    @{Row=19; B=0.00283870272538234; C=0.650055774033806},
    @{Row=20; B=0.00334363662860974; C=0.597313954514135},
    @{Row=21; B=0.00348133952728797; C=0.634096255422113},
    @{Row=22; B=0.00299047492118853; C=0.607996562858367},
    @{Row=23; B=0.00289508774243836; C=0.717438195996371},
    @{Row=24; B=0.00306790351104019; C=0.639330681291519},
    @{Row=25; B=0.00266115520863512; C=0.545079597601623},
    @{Row=26; B=0.030933790595386; C=0.716618950357581},
    @{Row=27; B=0.0328945591194969; C=0.674813101041297},
    @{Row=28; B=0.0281971644818462; C=0.728783640612074},
    @{Row=29; B=0.0535874616100596; C=0.581046802900461},
    @{Row=30; B=0.0334513327190158; C=0.75391440745423},
    @{Row=31; B=0.025156712332445; C=0.778394955788531},
    @{Row=32; B=0.0303442538926168; C=0.727860705106073},
    @{Row=33; B=0.0260115006088414; C=0.806954413652137},
    @{Row=34; B=0.022669706968705; C=0.79637964195887}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 2).Value = $u.B
    $ws.Cells.Item($u.Row, 3).Value = $u.C
}
